$d = $word.ActiveDocument

# Pull the canonical package XML (Content.WordOpenXML exposes the whole
# pkg:package, including /word/document.xml) so we can read each
# paragraph's authoritative OOXML without relying on the (unsupported)
# ParagraphFormat.ContextualSpacing property.
$full = $d.Content.WordOpenXML

$partMarker = '<pkg:part pkg:name="/word/document.xml"'
$partIdx = $full.IndexOf($partMarker)
$dataStart = $full.IndexOf("<pkg:xmlData>", $partIdx) + "<pkg:xmlData>".Length
$dataEnd = $full.IndexOf("</pkg:xmlData>", $dataStart)
$docXml = $full.Substring($dataStart, $dataEnd - $dataStart)

# Grab every top-level <w:p ...>...</w:p> paragraph block, in document order.
$paraMatches = [regex]::Matches($docXml, '<w:p(?:\s[^>]*)?(?:/>|>.*?</w:p>)', [System.Text.RegularExpressions.RegexOptions]::Singleline)

$count = $d.Paragraphs.Count
$i = 0
foreach ($m in $paraMatches) {
    if ($i -ge $count) { break }
    $i = $i + 1
    $frag = $m.Value

    if ($frag.IndexOf("contextualSpacing") -lt 0) { continue }

    # Strip the <w:contextualSpacing .../> element emitted by the GD export.
    $clean = $frag -replace '<w:contextualSpacing[^/]*/>', ''

    # Make the fragment namespace-self-contained before handing it to
    # InsertXML (it is spliced in on its own, outside the original
    # w:document element, so it needs its own xmlns declarations).
    $clean = $clean -replace '^<w:p ', '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" '

    $d.Paragraphs.Item($i).Range.InsertXML($clean)
}

Write-Output "paragraphs processed"
